$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 42.783101702509
$ws.Range("D11").Value = 269.322666666667

$ws.Range("C12").Value = 48.4471285842293
$ws.Range("D12").Value = 277.519333333333

$ws.Range("C13").Value = 0.0238963293650794
$ws.Range("D13").Value = 4.66966666666667

$ws.Range("C15").Value = 40.1516666666667
$ws.Range("D15").Value = 268.195333333333

$ws.Range("C16").Value = 13.4003042114696
$ws.Range("D16").Value = 291.598666666667
